# Update the two "Url" hyperlink cells (B2, B3) so their displayed text
# points at the new local auth endpoint instead of the old BodeWeb one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "http://localhost:3000/auth"
$ws.Range("B3").Value = "http://localhost:3000/auth"

# Re-normalize the cell formatting on the cells that were touched (B2, B3)
# plus the stray "Valor Esperado" numeric cell (F2) so they share the same
# plain bordered style as the rest of the data rows, instead of keeping
# their own now-redundant style entries.
$ws.Range("A2").Copy()
$ws.Range("B2").PasteSpecial(-4122)
$ws.Range("B3").PasteSpecial(-4122)
$ws.Range("F2").PasteSpecial(-4122)

$excel.CutCopyMode = 0
